$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.636.45"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "3.514.58"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.71"
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.07"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").Value = "3.509.70"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.17"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.23"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "4.081.46"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "605.88"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "3.515.20"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").Value = "70.745.27"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.72"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.878"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.10"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.49"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.10"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.45"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.00"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -5.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "621.38"
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0491"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.81"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0994"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -8.18%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.63"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "3.336.53"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "0.0₃0721"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.93"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.81"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  -6.75%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.96"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.01%  "
